# Applies the updated Betfair back/lay odds values (columns F:AO) for rows 2-7
# on Sheet1, per the source diff. Values are written directly as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 6).Value = 1.05  # F2
$ws.Cells.Item(2, 7).Value = 1.06  # G2
$ws.Cells.Item(2, 8).Value = 100  # H2
$ws.Cells.Item(2, 9).Value = 660  # I2
$ws.Cells.Item(2, 10).Value = 22  # J2
$ws.Cells.Item(2, 11).Value = 27  # K2
$ws.Cells.Item(2, 12).Value = 0  # L2
$ws.Cells.Item(2, 13).Value = 0  # M2
$ws.Cells.Item(2, 14).Value = 0  # N2
$ws.Cells.Item(2, 15).Value = 0  # O2
$ws.Cells.Item(2, 16).Value = 3.4  # P2
$ws.Cells.Item(2, 17).Value = 1.4  # Q2
$ws.Cells.Item(2, 18).Value = 1.54  # R2
$ws.Cells.Item(2, 19).Value = 2.74  # S2
$ws.Cells.Item(2, 20).Value = 2.28  # T2
$ws.Cells.Item(2, 21).Value = 1.48  # U2
$ws.Cells.Item(2, 22).Value = 1.01  # V2
$ws.Cells.Item(2, 23).Value = 18  # W2
$ws.Cells.Item(2, 24).Value = 1000  # X2
$ws.Cells.Item(2, 25).Value = 1000  # Y2
$ws.Cells.Item(2, 26).Value = 1000  # Z2
$ws.Cells.Item(2, 27).Value = 1000  # AA2
$ws.Cells.Item(2, 28).Value = 1000  # AB2
$ws.Cells.Item(2, 29).Value = 1000  # AC2
$ws.Cells.Item(2, 30).Value = 1000  # AD2
$ws.Cells.Item(2, 31).Value = 1000  # AE2
$ws.Cells.Item(2, 32).Value = 3.75  # AF2
$ws.Cells.Item(2, 33).Value = 7.4  # AG2
$ws.Cells.Item(2, 34).Value = 36  # AH2
$ws.Cells.Item(2, 35).Value = 230  # AI2
$ws.Cells.Item(2, 36).Value = 4.9  # AJ2
$ws.Cells.Item(2, 37).Value = 9.800000000000001  # AK2
$ws.Cells.Item(2, 38).Value = 55  # AL2
$ws.Cells.Item(2, 39).Value = 550  # AM2
$ws.Cells.Item(2, 40).Value = 6.6  # AN2
$ws.Cells.Item(2, 41).Value = 1000  # AO2
# Row 3
$ws.Cells.Item(3, 6).Value = 1.41  # F3
$ws.Cells.Item(3, 7).Value = 1.42  # G3
$ws.Cells.Item(3, 8).Value = 7.8  # H3
$ws.Cells.Item(3, 9).Value = 9.6  # I3
$ws.Cells.Item(3, 10).Value = 5.6  # J3
$ws.Cells.Item(3, 11).Value = 6  # K3
$ws.Cells.Item(3, 12).Value = 1.3  # L3
$ws.Cells.Item(3, 13).Value = 1.03  # M3
$ws.Cells.Item(3, 14).Value = 5.8  # N3
$ws.Cells.Item(3, 15).Value = 1.19  # O3
$ws.Cells.Item(3, 16).Value = 2.6  # P3
$ws.Cells.Item(3, 17).Value = 1.58  # Q3
$ws.Cells.Item(3, 18).Value = 1.6  # R3
$ws.Cells.Item(3, 19).Value = 2.54  # S3
$ws.Cells.Item(3, 20).Value = 1.81  # T3
$ws.Cells.Item(3, 21).Value = 2.16  # U3
$ws.Cells.Item(3, 23).Value = 3.35  # W3
$ws.Cells.Item(3, 24).Value = 34  # X3
$ws.Cells.Item(3, 25).Value = 40  # Y3
$ws.Cells.Item(3, 26).Value = 200  # Z3
$ws.Cells.Item(3, 27).Value = 300  # AA3
$ws.Cells.Item(3, 29).Value = 13.5  # AC3
$ws.Cells.Item(3, 31).Value = 110  # AE3
$ws.Cells.Item(3, 32).Value = 10  # AF3
$ws.Cells.Item(3, 33).Value = 9.800000000000001  # AG3
$ws.Cells.Item(3, 35).Value = 290  # AI3
$ws.Cells.Item(3, 36).Value = 12.5  # AJ3
$ws.Cells.Item(3, 37).Value = 13.5  # AK3
$ws.Cells.Item(3, 38).Value = 29  # AL3
$ws.Cells.Item(3, 39).Value = 990  # AM3
$ws.Cells.Item(3, 40).Value = 5.7  # AN3
$ws.Cells.Item(3, 41).Value = 160  # AO3
# Row 4
$ws.Cells.Item(4, 6).Value = 2.1  # F4
$ws.Cells.Item(4, 7).Value = 2.16  # G4
$ws.Cells.Item(4, 8).Value = 3.5  # H4
$ws.Cells.Item(4, 9).Value = 3.7  # I4
$ws.Cells.Item(4, 10).Value = 3.9  # J4
$ws.Cells.Item(4, 11).Value = 4  # K4
$ws.Cells.Item(4, 14).Value = 5  # N4
$ws.Cells.Item(4, 15).Value = 1.23  # O4
$ws.Cells.Item(4, 16).Value = 2.32  # P4
$ws.Cells.Item(4, 17).Value = 1.66  # Q4
$ws.Cells.Item(4, 18).Value = 1.54  # R4
$ws.Cells.Item(4, 19).Value = 2.66  # S4
$ws.Cells.Item(4, 20).Value = 1.65  # T4
$ws.Cells.Item(4, 21).Value = 2.44  # U4
$ws.Cells.Item(4, 22).Value = 1.37  # V4
$ws.Cells.Item(4, 23).Value = 1.86  # W4
$ws.Cells.Item(4, 24).Value = 22  # X4
$ws.Cells.Item(4, 26).Value = 40  # Z4
$ws.Cells.Item(4, 27).Value = 65  # AA4
$ws.Cells.Item(4, 29).Value = 9.6  # AC4
$ws.Cells.Item(4, 30).Value = 14.5  # AD4
$ws.Cells.Item(4, 31).Value = 36  # AE4
$ws.Cells.Item(4, 32).Value = 15.5  # AF4
$ws.Cells.Item(4, 34).Value = 16  # AH4
$ws.Cells.Item(4, 35).Value = 44  # AI4
$ws.Cells.Item(4, 36).Value = 27  # AJ4
$ws.Cells.Item(4, 38).Value = 28  # AL4
$ws.Cells.Item(4, 39).Value = 65  # AM4
$ws.Cells.Item(4, 40).Value = 11.5  # AN4
$ws.Cells.Item(4, 41).Value = 30  # AO4
# Row 5
$ws.Cells.Item(5, 6).Value = 2.14  # F5
$ws.Cells.Item(5, 7).Value = 2.22  # G5
$ws.Cells.Item(5, 8).Value = 4  # H5
$ws.Cells.Item(5, 9).Value = 4.5  # I5
$ws.Cells.Item(5, 10).Value = 3.2  # J5
$ws.Cells.Item(5, 11).Value = 3.4  # K5
$ws.Cells.Item(5, 12).Value = 1.52  # L5
$ws.Cells.Item(5, 14).Value = 2.88  # N5
$ws.Cells.Item(5, 15).Value = 1.5  # O5
$ws.Cells.Item(5, 16).Value = 1.61  # P5
$ws.Cells.Item(5, 17).Value = 2.42  # Q5
$ws.Cells.Item(5, 20).Value = 2.08  # T5
$ws.Cells.Item(5, 21).Value = 1.79  # U5
$ws.Cells.Item(5, 22).Value = 1.29  # V5
$ws.Cells.Item(5, 23).Value = 1.81  # W5
$ws.Cells.Item(5, 24).Value = 10  # X5
$ws.Cells.Item(5, 25).Value = 12.5  # Y5
$ws.Cells.Item(5, 26).Value = 970  # Z5
$ws.Cells.Item(5, 27).Value = 1000  # AA5
$ws.Cells.Item(5, 28).Value = 7.4  # AB5
$ws.Cells.Item(5, 29).Value = 10  # AC5
$ws.Cells.Item(5, 30).Value = 18  # AD5
$ws.Cells.Item(5, 31).Value = 70  # AE5
$ws.Cells.Item(5, 32).Value = 12  # AF5
$ws.Cells.Item(5, 33).Value = 11.5  # AG5
$ws.Cells.Item(5, 34).Value = 24  # AH5
$ws.Cells.Item(5, 35).Value = 230  # AI5
$ws.Cells.Item(5, 36).Value = 32  # AJ5
$ws.Cells.Item(5, 37).Value = 28  # AK5
$ws.Cells.Item(5, 39).Value = 200  # AM5
$ws.Cells.Item(5, 40).Value = 1000  # AN5
$ws.Cells.Item(5, 41).Value = 100  # AO5
# Row 6
$ws.Cells.Item(6, 6).Value = 2.06  # F6
$ws.Cells.Item(6, 7).Value = 2.12  # G6
$ws.Cells.Item(6, 8).Value = 4.5  # H6
$ws.Cells.Item(6, 9).Value = 5  # I6
$ws.Cells.Item(6, 10).Value = 3.15  # J6
$ws.Cells.Item(6, 11).Value = 3.35  # K6
$ws.Cells.Item(6, 14).Value = 2.52  # N6
$ws.Cells.Item(6, 15).Value = 1.58  # O6
$ws.Cells.Item(6, 16).Value = 1.49  # P6
$ws.Cells.Item(6, 17).Value = 2.78  # Q6
$ws.Cells.Item(6, 18).Value = 1.17  # R6
$ws.Cells.Item(6, 19).Value = 5.9  # S6
$ws.Cells.Item(6, 20).Value = 2.26  # T6
$ws.Cells.Item(6, 22).Value = 1.27  # V6
$ws.Cells.Item(6, 23).Value = 1.9  # W6
$ws.Cells.Item(6, 25).Value = 980  # Y6
$ws.Cells.Item(6, 28).Value = 11  # AB6
$ws.Cells.Item(6, 29).Value = 14  # AC6
$ws.Cells.Item(6, 32).Value = 21  # AF6
$ws.Cells.Item(6, 33).Value = 980  # AG6
$ws.Cells.Item(6, 36).Value = 95  # AJ6
# Row 7
$ws.Cells.Item(7, 6).Value = 3.65  # F7
$ws.Cells.Item(7, 8).Value = 2.2  # H7
$ws.Cells.Item(7, 9).Value = 2.44  # I7
$ws.Cells.Item(7, 10).Value = 3  # J7
$ws.Cells.Item(7, 11).Value = 3.4  # K7
$ws.Cells.Item(7, 12).Value = 1.55  # L7
$ws.Cells.Item(7, 13).Value = 1.11  # M7
$ws.Cells.Item(7, 14).Value = 2.92  # N7
$ws.Cells.Item(7, 15).Value = 1.46  # O7
$ws.Cells.Item(7, 16).Value = 1.61  # P7
$ws.Cells.Item(7, 17).Value = 2.44  # Q7
$ws.Cells.Item(7, 18).Value = 1.22  # R7
$ws.Cells.Item(7, 19).Value = 4.6  # S7
$ws.Cells.Item(7, 20).Value = 1.98  # T7
$ws.Cells.Item(7, 22).Value = 1.7  # V7
$ws.Cells.Item(7, 23).Value = 1.32  # W7
$ws.Cells.Item(7, 25).Value = 980  # Y7
$ws.Cells.Item(7, 26).Value = 1000  # Z7
$ws.Cells.Item(7, 29).Value = 1000  # AC7
$ws.Cells.Item(7, 30).Value = 12.5  # AD7
$ws.Cells.Item(7, 32).Value = 1000  # AF7
$ws.Cells.Item(7, 33).Value = 980  # AG7
